$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 17817
$ws.Range("E2").Value = 1159
$ws.Range("F2").Value = 1159
$ws.Range("G2").Value = 1251
$ws.Range("H2").Value = 941
$ws.Range("I2").Value = 938
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = 13516
$ws.Range("L2").Value = 5282
$ws.Range("M2").Value = 8233
$ws.Range("N2").Value = 8181
$ws.Range("O2").Value = 52
$ws.Range("P2").Value = 172
$ws.Range("Q2").Value = 1174
$ws.Range("R2").Value = -1517
$ws.Range("S2").Value = 341
$ws.Range("T2").Value = 978
$ws.Range("U2").Value = 196
$ws.Range("V2").Value = 1635
$ws.Range("W2").Value = 6.5
$ws.Range("X2").Value = 5.28
$ws.Range("Y2").Value = 12.06
$ws.Range("Z2").Value = 7.36
$ws.Range("AA2").Value = 64.16
$ws.Range("AB2").Value = 4562.69
$ws.Range("AC2").Value = 27258
$ws.Range("AD2").Value = 17.83
$ws.Range("AE2").Value = 242272
$ws.Range("AF2").Value = 2.01
$ws.Range("AG2").Value = 4000
$ws.Range("AH2").Value = 0.82
$ws.Range("AI2").Value = 14.41
$ws.Range("AJ2").Value = 3440000

# Row 3
$ws.Range("D3").Value = 18831
$ws.Range("E3").Value = 1334
$ws.Range("F3").Value = 1334
$ws.Range("G3").Value = 1429
$ws.Range("H3").Value = 1049
$ws.Range("I3").Value = 1045
$ws.Range("J3").Value = 5
$ws.Range("K3").Value = 14843
$ws.Range("L3").Value = 5352
$ws.Range("M3").Value = 9491
$ws.Range("N3").Value = 9434
$ws.Range("O3").Value = 57
$ws.Range("P3").Value = 172
$ws.Range("Q3").Value = 1115
$ws.Range("R3").Value = -814
$ws.Range("S3").Value = -372
$ws.Range("T3").Value = 812
$ws.Range("U3").Value = 303
$ws.Range("V3").Value = 1433
$ws.Range("W3").Value = 7.08
$ws.Range("X3").Value = 5.57
$ws.Range("Y3").Value = 11.86
$ws.Range("Z3").Value = 7.4
$ws.Range("AA3").Value = 56.39
$ws.Range("AB3").Value = 5049.33
$ws.Range("AC3").Value = 30365
$ws.Range("AD3").Value = 40.34
$ws.Range("AE3").Value = 279364
$ws.Range("AF3").Value = 4.38
$ws.Range("AG3").Value = 5200
$ws.Range("AH3").Value = 0.42
$ws.Range("AI3").Value = 16.81
$ws.Range("AJ3").Value = 3440000

# Row 4
$ws.Range("D4").Value = 20107
$ws.Range("E4").Value = 1425
$ws.Range("F4").Value = 1425
$ws.Range("G4").Value = 1834
$ws.Range("H4").Value = 1380
$ws.Range("I4").Value = 1375
$ws.Range("J4").Value = 5
$ws.Range("K4").Value = 15927
$ws.Range("L4").Value = 5575
$ws.Range("M4").Value = 10352
$ws.Range("N4").Value = 10292
$ws.Range("O4").Value = 60
$ws.Range("P4").Value = 172
$ws.Range("Q4").Value = 1489
$ws.Range("R4").Value = -1334
$ws.Range("S4").Value = -402
$ws.Range("T4").Value = 843
$ws.Range("U4").Value = 645
$ws.Range("V4").Value = 1359
$ws.Range("W4").Value = 7.09
$ws.Range("X4").Value = 6.86
$ws.Range("Y4").Value = 13.94
$ws.Range("Z4").Value = 8.970000000000001
$ws.Range("AA4").Value = 53.85
$ws.Range("AB4").Value = 5713.2
$ws.Range("AC4").Value = 39977
$ws.Range("AD4").Value = 16.56
$ws.Range("AE4").Value = 305644
$ws.Range("AF4").Value = 2.17
$ws.Range("AG4").Value = 6800
$ws.Range("AH4").Value = 1.03
$ws.Range("AI4").Value = 16.65
$ws.Range("AJ4").Value = 3440000

# Row 5
$ws.Range("D5").Value = 21262
$ws.Range("E5").Value = 1461
$ws.Range("F5").Value = 1461
$ws.Range("G5").Value = 1803
$ws.Range("H5").Value = 1324
$ws.Range("I5").Value = 1313
$ws.Range("J5").Value = 11
$ws.Range("K5").Value = 17217
$ws.Range("L5").Value = 5701
$ws.Range("M5").Value = 11516
$ws.Range("N5").Value = 10980
$ws.Range("O5").Value = 536
$ws.Range("P5").Value = 172
$ws.Range("Q5").Value = 1086
$ws.Range("R5").Value = -846
$ws.Range("S5").Value = -499
$ws.Range("T5").Value = 777
$ws.Range("U5").Value = 309
$ws.Range("V5").Value = 1200
$ws.Range("W5").Value = 6.87
$ws.Range("X5").Value = 6.23
$ws.Range("Y5").Value = 12.34
$ws.Range("Z5").Value = 7.99
$ws.Range("AA5").Value = 49.51
$ws.Range("AB5").Value = 6478.79
$ws.Range("AC5").Value = 38161
$ws.Range("AD5").Value = 21.12
$ws.Range("AE5").Value = 326055
$ws.Range("AF5").Value = 2.47
$ws.Range("AG5").Value = 7000
$ws.Range("AH5").Value = 0.87
$ws.Range("AI5").Value = 17.96
$ws.Range("AJ5").Value = 3440000

# Row 6
$ws.Range("D6").Value = 22468
$ws.Range("E6").Value = 1517
$ws.Range("F6").Value = 1517
$ws.Range("G6").Value = 2018
$ws.Range("H6").Value = 1608
$ws.Range("I6").Value = 1600
$ws.Range("K6").Value = 20635
$ws.Range("L6").Value = 7878
$ws.Range("M6").Value = 12757
$ws.Range("N6").Value = 12080
$ws.Range("P6").Value = 180
$ws.Range("Q6").Value = 1969
$ws.Range("R6").Value = -2719
$ws.Range("S6").Value = 699
$ws.Range("T6").Value = 2011
$ws.Range("U6").Value = -42
$ws.Range("V6").Value = 2711
$ws.Range("W6").Value = 6.75
$ws.Range("X6").Value = 7.16
$ws.Range("Y6").Value = 13.87
$ws.Range("Z6").Value = 8.5
$ws.Range("AA6").Value = 61.76
$ws.Range("AB6").Value = 7722.83
$ws.Range("AC6").Value = 45925
$ws.Range("AD6").Value = 15.76
$ws.Range("AE6").Value = 356097
$ws.Range("AF6").Value = 2.03
$ws.Range("AG6").Value = 7500
$ws.Range("AH6").Value = 1.04
$ws.Range("AI6").Value = 15.91
$ws.Range("AJ6").Value = 3605237

# Row 7
$ws.Range("D7").Value = 23508
$ws.Range("E7").Value = 1591
$ws.Range("G7").Value = 1601
$ws.Range("H7").Value = 1143
$ws.Range("I7").Value = 1135
$ws.Range("K7").Value = 21632
$ws.Range("L7").Value = 8070
$ws.Range("M7").Value = 13560
$ws.Range("N7").Value = 12861
$ws.Range("P7").Value = 180
$ws.Range("Q7").Value = 1387
$ws.Range("R7").Value = -744
$ws.Range("S7").Value = -423
$ws.Range("T7").Value = 1363
$ws.Range("U7").Value = 715
$ws.Range("W7").Value = 6.77
$ws.Range("X7").Value = 4.86
$ws.Range("Y7").Value = 9.1
$ws.Range("Z7").Value = 5.41
$ws.Range("AA7").Value = 59.51
$ws.Range("AC7").Value = 31482
$ws.Range("AD7").Value = 16.1
$ws.Range("AE7").Value = 379121
$ws.Range("AF7").Value = 1.34
$ws.Range("AG7").Value = 7625
$ws.Range("AH7").Value = 1.5
$ws.Range("AI7").Value = 24.22

# Row 8
$ws.Range("D8").Value = 24189
$ws.Range("E8").Value = 1675
$ws.Range("G8").Value = 1761
$ws.Range("H8").Value = 1321
$ws.Range("I8").Value = 1294
$ws.Range("K8").Value = 22644
$ws.Range("L8").Value = 8145
$ws.Range("M8").Value = 14498
$ws.Range("N8").Value = 13764
$ws.Range("P8").Value = 180
$ws.Range("Q8").Value = 1735
$ws.Range("R8").Value = -1021
$ws.Range("S8").Value = -383
$ws.Range("T8").Value = 1284
$ws.Range("U8").Value = 941
$ws.Range("W8").Value = 6.92
$ws.Range("X8").Value = 5.46
$ws.Range("Y8").Value = 9.720000000000001
$ws.Range("Z8").Value = 5.96
$ws.Range("AA8").Value = 56.18
$ws.Range("AC8").Value = 35892
$ws.Range("AD8").Value = 14.13
$ws.Range("AE8").Value = 405717
$ws.Range("AF8").Value = 1.25
$ws.Range("AG8").Value = 8000
$ws.Range("AH8").Value = 1.58
$ws.Range("AI8").Value = 22.29

# Row 9
$ws.Range("D9").Value = 24844
$ws.Range("E9").Value = 1762
$ws.Range("G9").Value = 1838
$ws.Range("H9").Value = 1371
$ws.Range("I9").Value = 1347
$ws.Range("K9").Value = 23815
$ws.Range("L9").Value = 8222
$ws.Range("M9").Value = 15596
$ws.Range("N9").Value = 14829
$ws.Range("P9").Value = 180
$ws.Range("Q9").Value = 1787
$ws.Range("R9").Value = -1002
$ws.Range("S9").Value = -411
$ws.Range("T9").Value = 1300
$ws.Range("U9").Value = 964
$ws.Range("W9").Value = 7.09
$ws.Range("X9").Value = 5.52
$ws.Range("Y9").Value = 9.42
$ws.Range("Z9").Value = 5.9
$ws.Range("AA9").Value = 52.72
$ws.Range("AC9").Value = 37362
$ws.Range("AD9").Value = 13.57
$ws.Range("AE9").Value = 437126
$ws.Range("AF9").Value = 1.16
$ws.Range("AG9").Value = 8250
$ws.Range("AH9").Value = 1.63
$ws.Range("AI9").Value = 22.08
